$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.Value = "'" + $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "58.325.55"
Set-TextValue "E2" "  -1.42%  "

Set-TextValue "D3" "2.569.97"
Set-TextValue "E3" "  -3.18%  "

Set-TextValue "E4" "  -0.04%  "

Set-TextValue "D5" "522.58"
Set-TextValue "E5" "  -0.32%  "

Set-TextValue "D6" "140.99"
Set-TextValue "E6" "  -2.70%  "

Set-TextValue "E7" "  -0.01%  "

Set-TextValue "D8" "0.557"
Set-TextValue "E8" "  -2.59%  "

Set-TextValue "D9" "6.57"
Set-TextValue "E9" "  -6.89%  "

Set-TextValue "D10" "0.0999"
Set-TextValue "E10" "  -2.79%  "

Set-TextValue "D11" "0.325"
Set-TextValue "E11" "  -2.94%  "

Set-TextValue "E12" "  +0.13%  "

Set-TextValue "D13" "3.030.47"
Set-TextValue "E13" "  -2.90%  "

Set-TextValue "D14" "58.303.31"
Set-TextValue "E14" "  -1.49%  "

Set-TextValue "D15" "20.04"
Set-TextValue "E15" "  -5.23%  "

Set-TextValue "D16" "2.595.13"
Set-TextValue "E16" "  -3.54%  "

Set-TextValue "D17" "0.0000132"
Set-TextValue "E17" "  -3.33%  "

Set-TextValue "D18" "335.48"
Set-TextValue "E18" "  -1.56%  "

Set-TextValue "D19" "4.30"
Set-TextValue "E19" "  -2.14%  "

Set-TextValue "D20" "10.14"
Set-TextValue "E20" "  -2.36%  "

Set-TextValue "D21" "6.14"
Set-TextValue "E21" "  -3.73%  "

Set-TextValue "E22" "  +0.21%  "

Set-TextValue "D23" "65.14"
Set-TextValue "E23" "  +1.12%  "

Set-TextValue "B24" "WrappedeETH"
Set-TextValue "C24" "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue "D24" "2.744.39"
Set-TextValue "E24" "  -1.08%  "

Set-TextValue "B25" "Polygon"
Set-TextValue "C25" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D25" "0.403"
Set-TextValue "E25" "  -2.80%  "

Set-TextValue "B26" "Kaspa"
Set-TextValue "C26" "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D26" "0.163"
Set-TextValue "E26" "  -1.99%  "

Set-TextValue "D27" "1.00"
Set-TextValue "E27" "  +0.37%  "

Set-TextValue "D28" "6.96"
Set-TextValue "E28" "  -2.17%  "

Set-TextValue "D29" "0.0₃0757"
Set-TextValue "E29" "  -5.95%  "

Set-TextValue "D30" "0.998"
Set-TextValue "E30" "  +0.01%  "

Set-TextValue "D31" "6.15"
Set-TextValue "E31" "  -8.08%  "

Set-TextValue "D32" "1.57"
Set-TextValue "E32" "  -1.18%  "

Set-TextValue "B33" "Monero"
Set-TextValue "C33" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D33" "149.41"
Set-TextValue "E33" "  +0.08%  "

Set-TextValue "B34" "EthereumClassic"
Set-TextValue "C34" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D34" "18.49"
Set-TextValue "E34" "  -1.92%  "

Set-TextValue "D35" "3.99"
Set-TextValue "E35" "  -4.33%  "

Set-TextValue "D36" "1.13"
Set-TextValue "E36" "  -5.43%  "

Set-TextValue "B37" "SuiNetwork"
Set-TextValue "C37" "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
Set-TextValue "D37" "0.830"
Set-TextValue "E37" "  -7.61%  "

Set-TextValue "B38" "OKB"
Set-TextValue "C38" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D38" "35.87"
Set-TextValue "E38" "  -2.27%  "

Set-TextValue "D39" "0.824"
Set-TextValue "E39" "  -7.11%  "

Set-TextValue "D40" "1.41"
Set-TextValue "E40" "  -4.74%  "

Set-TextValue "D41" "3.47"
Set-TextValue "E41" "  -3.41%  "

Set-TextValue "D42" "0.999"
Set-TextValue "E42" "  +0.01%  "

Set-TextValue "D43" "0.0955"
Set-TextValue "E43" "  -1.73%  "

Set-TextValue "B44" "WhiteBITCoin"
Set-TextValue "C44" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D44" "10.66"
Set-TextValue "E44" "  +1.25%  "

Set-TextValue "B45" "Mantle"
Set-TextValue "C45" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D45" "0.584"
Set-TextValue "E45" "  -5.57%  "

Set-TextValue "D46" "263.95"
Set-TextValue "E46" "  -4.11%  "

Set-TextValue "B47" "Hedera"
Set-TextValue "C47" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D47" "0.0518"
Set-TextValue "E47" "  -3.03%  "

Set-TextValue "B48" "EnergySwap"
Set-TextValue "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "18.57"
Set-TextValue "E48" "  -6.88%  "

Set-TextValue "D49" "1.984.61"
Set-TextValue "E49" "  -2.22%  "

Set-TextValue "D50" "0.0222"
Set-TextValue "E50" "  -3.32%  "

Set-TextValue "B51" "RenderToken"
Set-TextValue "C51" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D51" "4.51"
Set-TextValue "E51" "  -5.69%  "
